$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '65.570.02'
$ws.Range('E2').Value = '  -0.44%  '
$ws.Range('D3').Value = '2.944.69'
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '570.79'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -2.04%  '
$ws.Range('E6').Value = '  +0.57%  '
$ws.Range('E7').Value = '  -0.08%  '
$ws.Range('E8').Value = '  -0.56%  '
$ws.Range('D9').Value = '2.941.47'
$ws.Range('E9').Value = '  -2.11%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.69'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.48%  '
$ws.Range('E11').Value = '  -3.74%  '
$ws.Range('E12').Value = '  +0.82%  '
$ws.Range('E13').Value = '  -2.91%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '34.86'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.56%  '
$ws.Range('E15').Value = '  -0.72%  '
$ws.Range('D16').Value = '65.553.09'
$ws.Range('E16').Value = '  -0.47%  '
$ws.Range('D17').Value = '3.434.69'
$ws.Range('E17').Value = '  -2.10%  '
$ws.Range('E18').Value = '  +0.74%  '
$ws.Range('D19').Value = '2.947.99'
$ws.Range('E19').Value = '  -2.13%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '15.71'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +12.30%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '445.82'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E22').Value = '  +0.75%  '
$ws.Range('E23').Value = '  -1.35%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '82.07'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.37%  '
$ws.Range('E25').Value = '  -1.48%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '12.21'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -1.60%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.01'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -6.68%  '
$ws.Range('E28').Value = '  -0.02%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.47'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +5.04%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '8.10'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.64%  '
$ws.Range('E31').Value = '  -0.28%  '
$ws.Range('E32').Value = '  -3.17%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.115'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +4.38%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '27.22'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.89%  '
$ws.Range('E35').Value = '  -0.01%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.974'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.92%  '
$ws.Range('E37').Value = '  -1.91%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '46.54'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +5.90%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '49.27'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.90%  '
$ws.Range('B40').Value = 'Stacks'
$ws.Range('C40').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.97'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -8.82%  '
$ws.Range('B41').Value = 'TheGraph'
$ws.Range('C41').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.303'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -2.27%  '
$ws.Range('E42').Value = '  -1.01%  '
$ws.Range('E43').Value = '  -5.01%  '
$ws.Range('E44').Value = '  +0.58%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '384.13'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -2.23%  '
$ws.Range('E46').Value = '  -1.10%  '
$ws.Range('D47').Value = '2.675.77'
$ws.Range('E47').Value = '  -4.12%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '133.87'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.77%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '24.02'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.89%  '
$ws.Range('E51').Value = '  +1.37%  '
